$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) type for numeric-looking Price values so Excel
# does not auto-convert them to numbers (which would also silently drop
# significant trailing zeros, e.g. "0.310" -> 0.31).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto-price refresh.
$ws.Range("D2").Value = "35.150.29"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "1.827.00"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "232.79"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "42.77"
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("D9").Value = "0.310"
$ws.Range("E9").Value = "  +5.99%  "
$ws.Range("D10").Value = "0.0689"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "2.098.23"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.843.28"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "11.18"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "0.664"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "4.69"
$ws.Range("E16").Value = "  +5.90%  "
$ws.Range("D17").Value = "35.101.82"
$ws.Range("D18").Value = "70.02"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "0.0₃0792"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").Value = "239.69"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "11.83"
$ws.Range("E21").Value = "  +5.99%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "4.59"
$ws.Range("E23").Value = "  +11.51%  "
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("D25").Value = "171.52"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "7.80"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "17.56"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  +30.61%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "3.352.05"
$ws.Range("E31").Value = "  +37.96%  "
$ws.Range("D32").Value = "0.0554"
$ws.Range("E32").Value = "  +7.32%  "
$ws.Range("D33").Value = "3.91"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").Value = "93.43"
$ws.Range("E36").Value = "  +9.78%  "
$ws.Range("D37").Value = "0.680"
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "1.324.11"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").Value = "0.992"
$ws.Range("E42").Value = "  +4.95%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "14.91"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "6.22"
$ws.Range("E47").Value = "  +7.96%  "
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").Value = "2.009.19"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "100.75"
$ws.Range("E51").Value = "  -0.48%  "
